# CS107-161: Adds Excel validation
#
# Refresh the student roster: rows 2-6 get new e-mail addresses for the
# existing students, rows 7-10 become four brand new students (name +
# e-mail), and rows 11-12 are re-affirmed as-is. Each newly introduced
# e-mail address also gets its own cell-level hyperlink (mirroring the
# ones already present on C2/C3/C11/C12), in addition to the pre-existing
# C4:C12 range hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C (e-mail) for the five already-existing students ---
$ws.Range("C2").Value = "pasne.d@husky.neu.edu"
$ws.Range("C3").Value = "sood.s@husky.neu.edu"
$ws.Range("C4").Value = "shail@ccs.neu.edu"
$ws.Range("C5").Value = "dave.v@husky.neu.edu"
$ws.Range("C6").Value = "snow.j@husky.neu.edu"

# --- Column B (name) for the four newly added students ---
$ws.Range("B7").Value = "Danny"
$ws.Range("B8").Value = "Erica"
$ws.Range("B9").Value = "Flurry"
$ws.Range("B10").Value = "Gara"

# --- Column C (e-mail) for the newly added students, plus re-affirm the
#     last two (unchanged) rows ---
$ws.Range("C7").Value = "danny.d@husky.neu.edu"
$ws.Range("C8").Value = "sniper.e@husky.neu.edu"
$ws.Range("C9").Value = "majin.f@husky.neu.edu"
$ws.Range("C10").Value = "hawking.g@husky.neu.edu"
$ws.Range("C11").Value = "max@x.com"
$ws.Range("C12").Value = "kat@x.com"

# --- New per-cell hyperlinks for the updated/new e-mail addresses ---
# (C2, C3, C11, C12 already carry their own hyperlink; C5 continues to
# fall back to the existing C4:C12 range hyperlink.)
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:shail@ccs.neu.edu")
$ws.Range("C4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:snow.j@husky.neu.edu")
$ws.Range("C6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:danny.d@husky.neu.edu")
$ws.Range("C7").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:sniper.e@husky.neu.edu")
$ws.Range("C8").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:majin.f@husky.neu.edu")
$ws.Range("C9").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:hawking.g@husky.neu.edu")
$ws.Range("C10").Style = "Hyperlink"

# --- Leave the selection where the editor last left it ---
$ws.Range("C16").Select()
